$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 54.69462833333333
$ws.Range("H2").Value = 164.083885
$ws.Range("I2").Value = 0.2790924419198448
$ws.Range("J2").Value = 0.2790924419198448
$ws.Range("M2").Value = 27.85106533333333
$ws.Range("N2").Value = 83.553196
$ws.Range("O2").Value = 0.1861900221007236
$ws.Range("P2").Value = 0.1861900221007236
$ws.Range("Q2").Value = 1523.303667094051
$ws.Range("R2").Value = 13709.73300384646
$ws.Range("S2").Value = 0.05196422792920082
$ws.Range("T2").Value = 0.05196422792920082
$ws.Range("G3").Value = 54.69462833333333
$ws.Range("H3").Value = 164.083885
$ws.Range("I3").Value = 0.2790924419198448
$ws.Range("J3").Value = 0.2790924419198448
$ws.Range("O3").Value = 0.4727109026912454
$ws.Range("P3").Value = 0.4727109026912454
$ws.Range("Q3").Value = 3867.458864983479
$ws.Range("R3").Value = 34807.12978485131
$ws.Range("S3").Value = 0.1319300401542338
$ws.Range("T3").Value = 0.1319300401542338
$ws.Range("G4").Value = 54.69462833333333
$ws.Range("H4").Value = 164.083885
$ws.Range("I4").Value = 0.2790924419198448
$ws.Range("J4").Value = 0.2790924419198448
$ws.Range("M4").Value = 15.018964
$ws.Range("N4").Value = 45.056892
$ws.Range("O4").Value = 0.1004048213460311
$ws.Range("P4").Value = 0.1004048213460311
$ws.Range("Q4").Value = 821.4566539317133
$ws.Range("R4").Value = 7393.10988538542
$ws.Range("S4").Value = 0.02802222676998956
$ws.Range("T4").Value = 0.02802222676998956
$ws.Range("G5").Value = 54.69462833333333
$ws.Range("H5").Value = 164.083885
$ws.Range("I5").Value = 0.2790924419198448
$ws.Range("J5").Value = 0.2790924419198448
$ws.Range("M5").Value = 36.00403133333333
$ws.Range("N5").Value = 108.012094
$ws.Range("O5").Value = 0.2406942538619999
$ws.Range("P5").Value = 0.2406942538619999
$ws.Range("Q5").Value = 1969.227112278354
$ws.Range("R5").Value = 17723.04401050519
$ws.Range("S5").Value = 0.06717594706642059
$ws.Range("T5").Value = 0.06717594706642059
$ws.Range("G6").Value = 19.32115333333334
$ws.Range("H6").Value = 57.96346000000001
$ws.Range("I6").Value = 0.09859081282432611
$ws.Range("J6").Value = 0.09859081282432611
$ws.Range("M6").Value = 27.85106533333333
$ws.Range("N6").Value = 83.553196
$ws.Range("O6").Value = 0.1861900221007236
$ws.Range("P6").Value = 0.1861900221007236
$ws.Range("Q6").Value = 538.114703802018
$ws.Range("R6").Value = 4843.032334218161
$ws.Range("S6").Value = 0.01835662561868958
$ws.Range("T6").Value = 0.01835662561868958
$ws.Range("G7").Value = 19.32115333333334
$ws.Range("H7").Value = 57.96346000000001
$ws.Range("I7").Value = 0.09859081282432611
$ws.Range("J7").Value = 0.09859081282432611
$ws.Range("O7").Value = 0.4727109026912454
$ws.Range("P7").Value = 0.4727109026912454
$ws.Range("Q7").Value = 1366.199351155754
$ws.Range("S7").Value = 0.04660495212725081
$ws.Range("T7").Value = 0.04660495212725081
$ws.Range("G8").Value = 19.32115333333334
$ws.Range("H8").Value = 57.96346000000001
$ws.Range("I8").Value = 0.09859081282432611
$ws.Range("J8").Value = 0.09859081282432611
$ws.Range("M8").Value = 15.018964
$ws.Range("N8").Value = 45.056892
$ws.Range("O8").Value = 0.1004048213460311
$ws.Range("P8").Value = 0.1004048213460311
$ws.Range("Q8").Value = 290.1837063518134
$ws.Range("R8").Value = 2611.65335716632
$ws.Range("S8").Value = 0.00989899294798645
$ws.Range("T8").Value = 0.00989899294798645
$ws.Range("G9").Value = 19.32115333333334
$ws.Range("H9").Value = 57.96346000000001
$ws.Range("I9").Value = 0.09859081282432611
$ws.Range("J9").Value = 0.09859081282432611
$ws.Range("M9").Value = 36.00403133333333
$ws.Range("N9").Value = 108.012094
$ws.Range("O9").Value = 0.2406942538619999
$ws.Range("P9").Value = 0.2406942538619999
$ws.Range("Q9").Value = 695.6394100094712
$ws.Range("R9").Value = 6260.75469008524
$ws.Range("S9").Value = 0.02373024213039927
$ws.Range("T9").Value = 0.02373024213039927
$ws.Range("G10").Value = 11.023718
$ws.Range("H10").Value = 33.071154
$ws.Range("I10").Value = 0.05625116157486912
$ws.Range("J10").Value = 0.05625116157486911
$ws.Range("M10").Value = 27.85106533333333
$ws.Range("N10").Value = 83.553196
$ws.Range("O10").Value = 0.1861900221007236
$ws.Range("P10").Value = 0.1861900221007236
$ws.Range("Q10").Value = 307.0222902342427
$ws.Range("R10").Value = 2763.200612108184
$ws.Range("S10").Value = 0.01047340501681625
$ws.Range("T10").Value = 0.01047340501681625
$ws.Range("G11").Value = 11.023718
$ws.Range("H11").Value = 33.071154
$ws.Range("I11").Value = 0.05625116157486912
$ws.Range("J11").Value = 0.05625116157486911
$ws.Range("O11").Value = 0.4727109026912454
$ws.Range("P11").Value = 0.4727109026912454
$ws.Range("Q11").Value = 779.4874415152581
$ws.Range("R11").Value = 7015.386973637322
$ws.Range("S11").Value = 0.02659053736548748
$ws.Range("T11").Value = 0.02659053736548747
$ws.Range("G12").Value = 11.023718
$ws.Range("H12").Value = 33.071154
$ws.Range("I12").Value = 0.05625116157486912
$ws.Range("J12").Value = 0.05625116157486911
$ws.Range("M12").Value = 15.018964
$ws.Range("N12").Value = 45.056892
$ws.Range("O12").Value = 0.1004048213460311
$ws.Range("P12").Value = 0.1004048213460311
$ws.Range("Q12").Value = 165.564823788152
$ws.Range("R12").Value = 1490.083414093368
$ws.Range("S12").Value = 0.00564788782843146
$ws.Range("T12").Value = 0.00564788782843146
$ws.Range("G13").Value = 11.023718
$ws.Range("H13").Value = 33.071154
$ws.Range("I13").Value = 0.05625116157486912
$ws.Range("J13").Value = 0.05625116157486911
$ws.Range("M13").Value = 36.00403133333333
$ws.Range("N13").Value = 108.012094
$ws.Range("O13").Value = 0.2406942538619999
$ws.Range("P13").Value = 0.2406942538619999
$ws.Range("Q13").Value = 396.8982882818307
$ws.Range("R13").Value = 3572.084594536476
$ws.Range("S13").Value = 0.01353933136413392
$ws.Range("T13").Value = 0.01353933136413392
$ws.Range("G14").Value = 110.9336623333333
$ws.Range("H14").Value = 332.800987
$ws.Range("I14").Value = 0.5660655836809599
$ws.Range("J14").Value = 0.5660655836809599
$ws.Range("M14").Value = 27.85106533333333
$ws.Range("N14").Value = 83.553196
$ws.Range("O14").Value = 0.1861900221007236
$ws.Range("P14").Value = 0.1861900221007236
$ws.Range("Q14").Value = 3089.620677311606
$ws.Range("R14").Value = 27806.58609580445
$ws.Range("S14").Value = 0.1053957635360169
$ws.Range("T14").Value = 0.1053957635360169
$ws.Range("G15").Value = 110.9336623333333
$ws.Range("H15").Value = 332.800987
$ws.Range("I15").Value = 0.5660655836809599
$ws.Range("J15").Value = 0.5660655836809599
$ws.Range("O15").Value = 0.4727109026912454
$ws.Range("P15").Value = 0.4727109026912454
$ws.Range("Q15").Value = 7844.122702533531
$ws.Range("R15").Value = 70597.10432280178
$ws.Range("S15").Value = 0.2675853730442733
$ws.Range("T15").Value = 0.2675853730442733
$ws.Range("G16").Value = 110.9336623333333
$ws.Range("H16").Value = 332.800987
$ws.Range("I16").Value = 0.5660655836809599
$ws.Range("J16").Value = 0.5660655836809599
$ws.Range("M16").Value = 15.018964
$ws.Range("N16").Value = 45.056892
$ws.Range("O16").Value = 0.1004048213460311
$ws.Range("P16").Value = 0.1004048213460311
$ws.Range("Q16").Value = 1666.108680972489
$ws.Range("R16").Value = 14994.9781287524
$ws.Range("S16").Value = 0.05683571379962358
$ws.Range("T16").Value = 0.05683571379962358
$ws.Range("G17").Value = 110.9336623333333
$ws.Range("H17").Value = 332.800987
$ws.Range("I17").Value = 0.5660655836809599
$ws.Range("J17").Value = 0.5660655836809599
$ws.Range("M17").Value = 36.00403133333333
$ws.Range("N17").Value = 108.012094
$ws.Range("O17").Value = 0.2406942538619999
$ws.Range("P17").Value = 0.2406942538619999
$ws.Range("Q17").Value = 3994.059054570752
$ws.Range("R17").Value = 35946.53149113677
$ws.Range("S17").Value = 0.1362487333010461
$ws.Range("T17").Value = 0.1362487333010461
